$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the description for Merchant_Id (row 2) to clarify merchant id usage
$ws.Range("B2").Value = "Merchant ID, but could be of separate stores, use Merchant_name as unique Merchant identifier"

# Update the description for Merchant_Name (row 8) to clarify it should be used as the unique business identifier
$ws.Range("B8").Value = "Generic name of merchant, should be used as unique business identifier."

# Move the active selection, matching the recorded view state change
$ws.Range("G22").Select()
